$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) and Volume(1h) (E) updates ---
$ws.Range("D2").Value = "67.736.77"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "3.263.12"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'579.84"
$ws.Range("D6").Value = "'184.85"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  -4.31%  "
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("D11").Value = "'0.407"
$ws.Range("E11").Value = "  -3.87%  "
$ws.Range("D12").Value = "3.829.25"
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "'27.43"
$ws.Range("E14").Value = "  -6.04%  "
$ws.Range("D15").Value = "67.778.89"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("E16").Value = "  -3.15%  "
$ws.Range("D17").Value = "3.253.38"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("D18").Value = "'5.71"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("D19").Value = "'13.46"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("D20").Value = "'395.85"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "'7.56"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "'70.83"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("D24").Value = "'0.508"
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("E25").Value = "  -4.45%  "
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("D27").Value = "'9.48"
$ws.Range("E27").Value = "  -2.91%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("D30").Value = "'22.57"
$ws.Range("E30").Value = "  -2.50%  "
$ws.Range("D31").Value = "'5.46"
$ws.Range("E31").Value = "  -5.95%  "
$ws.Range("D32").Value = "'6.91"
$ws.Range("E32").Value = "  -3.81%  "
$ws.Range("D33").Value = "'0.998"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").Value = "'1.24"
$ws.Range("E34").Value = "  -5.84%  "
$ws.Range("D35").Value = "'163.36"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "'1.45"
$ws.Range("E36").Value = "  -5.96%  "
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("D38").Value = "'26.78"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("D39").Value = "'0.805"
$ws.Range("E39").Value = "  -3.92%  "
$ws.Range("D40").Value = "'4.50"
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("E41").Value = "  -4.89%  "
$ws.Range("D42").Value = "2.661.21"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("E45").Value = "  -7.94%  "
$ws.Range("D46").Value = "'24.67"
$ws.Range("E46").Value = "  -3.00%  "
$ws.Range("D47").Value = "'334.38"
$ws.Range("E47").Value = "  -2.35%  "
$ws.Range("E48").Value = "  -3.88%  "
$ws.Range("D49").Value = "'6.31"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").Value = "'0.970"
$ws.Range("E51").Value = "  -2.92%  "

# --- Row 43/44 swap (Hedera <-> OKB) ---
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "'40.67"
$ws.Range("E43").Value = "  -2.35%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "'0.0680"
$ws.Range("E44").Value = "  -2.05%  "
